$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the style used by
# the existing header row (H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data for the new I (I0) and J (IF) columns, rows 2-12.
$values = @(
    @(6, 6),
    @(6, 6),
    @(7, 7),
    @(5, 5),
    @(3, 4),
    @(7, 7),
    @(9, 9),
    @(5, 5),
    @(9, 9),
    @(6, 6),
    @(6, 6)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
